$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("B2").Value = 0.07881312460784126
$ws.Range("C2").Value = 0.07868058796995309
$ws.Range("B3").Value = 23.13921682692138
$ws.Range("C3").Value = 23.13993403873709
$ws.Range("B4").Value = 151.7001364441091
$ws.Range("C4").Value = 151.7000501335764
$ws.Range("B5").Value = 0.1858915563524067
$ws.Range("C5").Value = 0.1860573123609745
$ws.Range("B6").Value = 2.550002061758742
$ws.Range("C6").Value = 2.549960751924217
$ws.Range("B7").Value = 0.7718942589664032
$ws.Range("C7").Value = 0.7619561242120835
$ws.Range("B9").Value = 1.827249818929199
$ws.Range("C9").Value = 1.832533843157625
$ws.Range("B11").Value = 917.7353696920507
$ws.Range("C11").Value = 916.3224302482407
$ws.Range("B12").Value = 0.9923146550803537
$ws.Range("C12").Value = 0.99013189318705
$ws.Range("B13").Value = 1.052976498902757
$ws.Range("C13").Value = 0.9580825004213104
$ws.Range("B14").Value = 2.438337820041923
$ws.Range("C14").Value = 2.438558578626186
$ws.Range("B15").Value = 0.9145347591184037
$ws.Range("C15").Value = 0.9142978354870943
$ws.Range("B16").Value = 0.07859319507629441
$ws.Range("C16").Value = 0.07861429794597517
$ws.Range("B19").Value = 0.4584993047979745
$ws.Range("C19").Value = 0.4519771706760511
$ws.Range("B20").Value = 0.1304844436502391
$ws.Range("C20").Value = 0.1232399248041299
$ws.Range("B21").Value = 0.1291677274760185
$ws.Range("C21").Value = 0.1263881077842832
$ws.Range("B22").Value = 5.104926529470156
$ws.Range("C22").Value = 5.10987551632889
$ws.Range("B23").Value = -0.0242139505116133
$ws.Range("C23").Value = 0.00503156547321891
$ws.Range("B24").Value = 0.4089624395380318
$ws.Range("C24").Value = 0.4042404135919265
$ws.Range("B25").Value = 28.30648952303393
$ws.Range("C25").Value = 28.30598363412188
$ws.Range("B26").Value = 29.29966712599418
$ws.Range("C26").Value = 29.29963692361504
$ws.Range("B27").Value = 0.08208752150269018
$ws.Range("C27").Value = 0.08204491330088279
$ws.Range("B28").Value = 0.4543619591828427
$ws.Range("C28").Value = 0.4543825076948632
$ws.Range("B29").Value = 1.707315648853719
$ws.Range("C29").Value = 1.706783468221732
$ws.Range("B30").Value = 3.597811388646879
$ws.Range("C30").Value = 3.596831469393622
$ws.Range("B31").Value = 12.06603796442994
$ws.Range("C31").Value = 12.06554974753924
$ws.Range("B32").Value = 33.41805094524368
$ws.Range("C32").Value = 33.13508559218115
$ws.Range("B33").Value = 73866.9865992866
$ws.Range("C33").Value = 73852.06545373233
$ws.Range("B34").Value = 9.003243897527003
$ws.Range("C34").Value = 9.007231799129405
$ws.Range("B35").Value = 88.06246075798238
$ws.Range("C35").Value = 87.88212135037369
$ws.Range("B36").Value = 142.7725302575142
$ws.Range("C36").Value = 90.42901659229462